# Weekly update: a new price record is inserted at row 145 for
# "Terminal La Palmera de La Serena - Albahaca", pushing the existing
# rows 145:162 down to 146:163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 145 (shifts 145:162 -> 146:163, carries
# formatting from the row above, e.g. the date style on column D).
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new weekly record.
$ws.Cells.Item(145, 1).Value = 8
$ws.Cells.Item(145, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(145, 3).Value = "Coquimbo"
$ws.Cells.Item(145, 4).Value = 44984
$ws.Cells.Item(145, 5).Value = 4
$ws.Cells.Item(145, 6).Value = 100112052
$ws.Cells.Item(145, 7).Value = "Albahaca"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 800
$ws.Cells.Item(145, 11).Value = 4500
$ws.Cells.Item(145, 12).Value = 5000
$ws.Cells.Item(145, 13).Value = 4750
$ws.Cells.Item(145, 14).Value = "`$/docena de matas"
$ws.Cells.Item(145, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(145, 16).Value = 792
$ws.Cells.Item(145, 17).Value = 6
$ws.Cells.Item(145, 18).Value = "Hortaliza"
